$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '20.241.46'
$ws.Range("E2").Value = '  +1.11%  '

# Row 3
$ws.Range("D3").Value = '1.440.63'
$ws.Range("E3").Value = '  +1.44%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.008'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.81%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.9251'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -7.13%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '273.46'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3646'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.04%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3062'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.27%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '39.60'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.83%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.017'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.57%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.06498'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.18%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.9986'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.22%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.337'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.91%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '17.42'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.94%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.040'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.64%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.00001009'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.86%  '

# Row 17
$ws.Range("D17").Value = '1.438.64'
$ws.Range("E17").Value = '  +1.25%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.9435'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -5.34%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.05653'
$c.Style = "Normal"

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '68.42'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.07%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.359'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -4.70%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '14.19'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -4.02%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.77'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.54%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.245'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.20%  '

# Row 25
$ws.Range("D25").Value = '20.256.87'
$ws.Range("E25").Value = '  +0.83%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '140.32'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +3.36%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.025'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -11.03%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '16.88'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.84%  '

# Row 29
$ws.Range("D29").Value = '1.590.63'
$ws.Range("E29").Value = '  +0.73%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '110.22'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.12%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.018'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.29%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.781'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -11.25%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.07665'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.77%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.7734'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -7.76%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.451'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.76%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.05652'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -4.50%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.622'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -5.53%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.114'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.44%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01986'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -4.30%  '

# Row 40
$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.9362'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -6.00%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '10.13'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -5.43%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.1830'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -4.64%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '6.954'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -17.66%  '

# Row 44
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.5174'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.70%  '

# Row 45
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.469'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.64%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '11.66'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -5.71%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '114.15'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.42%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.5061'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.27%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.722'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.83%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06356'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.89%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.9892'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.87%  '
